$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: a new date column (AD = "25_05_2021") is appended after
# the previous last column (AC = "18_05_2021"), with one new death-count
# figure per age-group row plus a new "I alt" (total) sum row.

$ws.Range("AD1").Value = "25_05_2021"

$ws.Range("AD2").Value  = 1
$ws.Range("AD3").Value  = 0
$ws.Range("AD4").Value  = 0
$ws.Range("AD5").Value  = 7
$ws.Range("AD6").Value  = 8
$ws.Range("AD7").Value  = 64
$ws.Range("AD8").Value  = 210
$ws.Range("AD9").Value  = 666
$ws.Range("AD10").Value = 992
$ws.Range("AD11").Value = 562

$ws.Range("AD12").Formula = "=SUM(AD2:AD11)"

# Refresh the previous total-row cell's formula too (matches the source
# workbook, where it now stores its own literal SUM rather than sharing
# the formula group with the other weekly total cells).
$ws.Range("AC12").Formula = "=SUM(AC2:AC11)"

# Move the active selection to match the author's saved view state.
$ws.Range("AD17").Select()
